$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook is a localization-status "handback" report. A handoff run
# completed: the old source doc (dfed2ff0-...md) was re-handed-off under a
# fresh id (9e3cc3a9-...md) with new target .xlf artifacts + timestamps, and
# the stale "d482a737-...md / Handoff failed" row is gone entirely (its
# immediate successor, ".localization-config", shifts up one row).
# ---------------------------------------------------------------------------

$oldUuid = "dfed2ff0-0c8a-4795-8f89-04c3a8a66ec0"
$newUuid = "9e3cc3a9-430b-4a7b-98f3-43d1e4d708bc"
$oldHash = "94c5fc70554ab16089ed18696ace32ea269504c4"
$newHash = "8cabcc3d17f9b8490167c021781d42154f59c8ff"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0f7b00f8831d32edf94999a3440b039c8a0104e3/e2e/$newUuid.md"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0f7b00f8831d32edf94999a3440b039c8a0104e3/.localization-config"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5691fc7330a908c8a1ff397c05138848222e3ab5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$newUuid.$newHash.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ca30cdeab88bf5e67a5328bebceb807844f8128/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$newUuid.$newHash.de-de.xlf"

# Matches the workbook's custom "HyperLink" cell style (underline, RGB 6495ED).
$hlColor = 15570276

function Restyle-Hyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hlColor
}

# =============================== Sheet 1: Overview ==========================
$ws1 = $wb.Worksheets.Item(1)

# Hyperlink refs don't reflow when rows move, so drop them all up front and
# re-add clean ones (pointing at the post-shift rows) at the end.
$ws1.Hyperlinks.Delete()

# The "d482a737.../Handoff failed" row (row 3) is gone; row 4 shifts to row 3.
$ws1.Rows.Item(3).Delete()

$ws1.Range("A2").Value = "$newUuid.md"
$ws1.Range("A3").Value = ".localization-config"

$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdUrl, "", "", "$newUuid.md")
Restyle-Hyperlink $ws1.Range("A2")
$ws1.Hyperlinks.Add($ws1.Range("A3"), $cfgUrl, "", "", ".localization-config")
Restyle-Hyperlink $ws1.Range("A3")

# =============================== Sheet 2: zh-cn ==============================
$ws2 = $wb.Worksheets.Item(2)
$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(3).Delete()

$ws2.Range("A2").Value = "$newUuid.md"
$ws2.Range("C2").Value = "$newUuid.$newHash.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-08 12:07:55"
$ws2.Range("A3").Value = ".localization-config"

$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdUrl, "", "", "$newUuid.md")
Restyle-Hyperlink $ws2.Range("A2")
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zhXlfUrl, "", "", "$newUuid.$newHash.zh-cn.xlf")
Restyle-Hyperlink $ws2.Range("C2")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $cfgUrl, "", "", ".localization-config")
Restyle-Hyperlink $ws2.Range("A3")

# =============================== Sheet 3: de-de ==============================
$ws3 = $wb.Worksheets.Item(3)
$ws3.Hyperlinks.Delete()
$ws3.Rows.Item(3).Delete()

$ws3.Range("A2").Value = "$newUuid.md"
$ws3.Range("C2").Value = "$newUuid.$newHash.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-08 12:08:09"
$ws3.Range("A3").Value = ".localization-config"

$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdUrl, "", "", "$newUuid.md")
Restyle-Hyperlink $ws3.Range("A2")
$ws3.Hyperlinks.Add($ws3.Range("C2"), $deXlfUrl, "", "", "$newUuid.$newHash.de-de.xlf")
Restyle-Hyperlink $ws3.Range("C2")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $cfgUrl, "", "", ".localization-config")
Restyle-Hyperlink $ws3.Range("A3")
